# Update the table style used by the table on "slide 16" (Component 3
# summary slide) from the deck's custom table style to the built-in
# PowerPoint table style {B216833F-23D2-4EF1-9D0E-23D9FD228276}.
#
# We walk every slide/shape rather than hard-coding indices so the
# script still finds the table if shape ordering ever shifts.

$p = $ppt.ActivePresentation

$oldStyleId = "{8B964DD7-284D-469F-B8CA-A03F4789D4EF}"
$newStyleId = "{B216833F-23D2-4EF1-9D0E-23D9FD228276}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable -eq -1) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
